$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$words = @(
    'werden',
    'enden',
    'stehlen',
    'drehen',
    'trotzen',
    'sichern',
    'formen',
    'spielen',
    'runden',
    'heilen',
    'ärgern',
    'klappen',
    'hauen',
    'wenden',
    'fangen',
    'seufzen',
    'segnen',
    'wärmen',
    'wohnen',
    'trösten',
    'warnen',
    'doppeln',
    'haben',
    'schwören',
    'flüchten',
    'schnellen',
    'schwächen',
    'schmecken',
    'flehen',
    'fließen',
    'folgen',
    'scheitern',
    'starten',
    'wählen',
    'zögern',
    'lösen',
    'trennen',
    'sorgen',
    'dürfen',
    'wüten',
    'schwingen',
    'stecken',
    'hoffen',
    'töten',
    'klingeln',
    'freuen',
    'zielen',
    'orten',
    'wandern',
    'reizen',
    'jubeln',
    'zünden',
    'tropfen',
    'garen',
    'reisen',
    'heben',
    'weichen',
    'gelten',
    'prügeln',
    'bauen',
    'brauchen',
    'helfen',
    'dauern',
    'schreiten',
    'pflanzen',
    'landen',
    'wachsen',
    'filmen',
    'stammen',
    'bleiben',
    'fahren',
    'trauen',
    'wehtun',
    'wundern',
    'schultern',
    'ehren',
    'schenken',
    'streichen',
    'fallen',
    'schließen',
    'machen',
    'ahnen',
    'liegen',
    'geben',
    'warten',
    'treffen',
    'tollen',
    'betteln',
    'loben',
    'hupen',
    'bergen',
    'biegen',
    'quälen',
    'wagen',
    'jagen',
    'fischen',
    'locken',
    'segeln',
    'sterben',
    'schlucken',
    'räumen',
    'altern',
    'opfern',
    'leeren',
    'foltern',
    'schleppen',
    'lesen',
    'siegen',
    'meistern',
    'schrecken',
    'kichern',
    'boxen',
    'tauschen',
    'spinnen',
    'rufen',
    'klingen',
    'schlafen',
    'kümmern',
    'greifen',
    'graben',
    'mauern',
    'wiegen',
    'tanzen',
    'werfen',
    'bitten',
    'dringen',
    'grüßen',
    'lockern',
    'backen',
    'schalten',
    'schütteln',
    'pfeifen',
    'bellen',
    'gründen',
    'sinken',
    'rechnen',
    'hören',
    'achten',
    'suchen',
    'erben',
    'wirken',
    'stören',
    'sperren',
    'herrschen',
    'fällen',
    'kennen',
    'kämpfen',
    'planen',
    'scheinen',
    'flüstern',
    'wellen',
    'sprengen',
    'führen',
    'zählen',
    'platzen',
    'ändern',
    'schaden',
    'bluten',
    'münzen',
    'kosten',
    'stärken',
    'kehren',
    'irren',
    'rauchen',
    'klettern',
    'brauen',
    'malen',
    'rasen',
    'wehren',
    'saufen',
    'äußern',
    'schicken',
    'zahlen',
    'treiben',
    'mögen',
    'heulen',
    'weigern',
    'laufen',
    'arten',
    'zeigen',
    'ruhen',
    'decken',
    'feiern',
    'dienen',
    'stoppen',
    'spüren',
    'schulden',
    'knarren',
    'wechseln',
    'schämen',
    'liefern',
    'lügen'
)

for ($i = 0; $i -lt $words.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $words[$i]
}

